$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while guaranteeing it stays stored
# as text (not auto-converted to a number) and without leaving behind any
# permanent number-format / style change on the cell.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Updated coin price / 1h-volume figures (row order matches the sheet).
Set-TextValue 'D2' '26.076.51'
Set-TextValue 'E2' '  -0.62%  '
Set-TextValue 'D3' '1.659.31'
Set-TextValue 'E3' '  -1.29%  '
Set-TextValue 'E4' '  -0.15%  '
Set-TextValue 'D5' '207.85'
Set-TextValue 'E5' '  -1.24%  '
Set-TextValue 'D6' '0.5162'
Set-TextValue 'E6' '  -2.54%  '
Set-TextValue 'E7' '  -0.11%  '
Set-TextValue 'D8' '0.2579'
Set-TextValue 'E8' '  -3.89%  '
Set-TextValue 'D9' '0.06274'
Set-TextValue 'E9' '  -0.52%  '
Set-TextValue 'D10' '20.88'
Set-TextValue 'E10' '  -2.42%  '
Set-TextValue 'D11' '0.07519'
Set-TextValue 'E11' '  -0.27%  '
Set-TextValue 'D12' '1.667.65'
Set-TextValue 'E12' '  -0.87%  '
Set-TextValue 'D13' '4.387'
Set-TextValue 'E13' '  -1.98%  '
Set-TextValue 'D14' '0.5381'
Set-TextValue 'E14' '  -5.24%  '
Set-TextValue 'D15' '65.99'
Set-TextValue 'E15' '  -1.05%  '
Set-TextValue 'D16' '0.0₅7889'
Set-TextValue 'E16' '  -3.11%  '
Set-TextValue 'D17' '26.094.82'
Set-TextValue 'E17' '  -0.63%  '
Set-TextValue 'E18' '  -0.13%  '
Set-TextValue 'D19' '4.677'
Set-TextValue 'E19' '  -3.70%  '
Set-TextValue 'D20' '187.12'
Set-TextValue 'E20' '  -0.87%  '
Set-TextValue 'D21' '10.16'
Set-TextValue 'E21' '  -3.64%  '
Set-TextValue 'D22' '6.174'
Set-TextValue 'E22' '  -0.79%  '
Set-TextValue 'E23' '  -0.14%  '
Set-TextValue 'D24' '148.08'
Set-TextValue 'E24' '  +0.66%  '
Set-TextValue 'D25' '0.1205'
Set-TextValue 'E25' '  -4.42%  '
Set-TextValue 'D26' '7.357'
Set-TextValue 'E26' '  -3.83%  '
Set-TextValue 'E27' '  -1.97%  '
Set-TextValue 'D28' '1.379'
Set-TextValue 'E28' '  +2.18%  '
Set-TextValue 'D29' '0.06084'
Set-TextValue 'E29' '  -5.20%  '
Set-TextValue 'D30' '1.261'
Set-TextValue 'E30' '  -1.81%  '
Set-TextValue 'D31' '3.458'
Set-TextValue 'E31' '  -2.45%  '
Set-TextValue 'D32' '3.389'
Set-TextValue 'E32' '  -2.92%  '
Set-TextValue 'E33' '  -2.09%  '
Set-TextValue 'D34' '0.9796'
Set-TextValue 'E34' '  -3.53%  '
Set-TextValue 'E35' '  -1.07%  '
Set-TextValue 'E36' '  +1.21%  '
Set-TextValue 'D37' '0.5848'
Set-TextValue 'E37' '  -4.16%  '
Set-TextValue 'D38' '1.101.33'
Set-TextValue 'E38' '  -0.10%  '
Set-TextValue 'D39' '0.01590'
Set-TextValue 'E39' '  -1.63%  '
Set-TextValue 'D40' '5.956'
Set-TextValue 'E40' '  -3.56%  '
Set-TextValue 'D41' '0.8470'
Set-TextValue 'E41' '  -2.69%  '
Set-TextValue 'E42' '  -0.42%  '
Set-TextValue 'D43' '99.82'
Set-TextValue 'E43' '  -0.31%  '
Set-TextValue 'D44' '1.812.37'
Set-TextValue 'E44' '  -1.10%  '
Set-TextValue 'D45' '0.0₈108'
Set-TextValue 'E45' '  -1.24%  '
Set-TextValue 'B46' 'Aave'
Set-TextValue 'C46' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D46' '54.82'
Set-TextValue 'E46' '  -3.82%  '
Set-TextValue 'B47' 'Frax'
Set-TextValue 'C47' 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D47' '0.9858'
Set-TextValue 'E47' '  -1.53%  '
Set-TextValue 'D48' '7.962'
Set-TextValue 'E48' '  -0.59%  '
Set-TextValue 'D49' '0.05225'
Set-TextValue 'E49' '  -0.73%  '
Set-TextValue 'D50' '0.4239'
Set-TextValue 'E50' '  -0.72%  '
Set-TextValue 'D51' '5.846'
Set-TextValue 'E51' '  -2.22%  '
